# Replied to more posts.
# Fill in "Actual time length to complete" (column C) for two more
# discussion-post rows on the week3 sheet, then leave the selection
# where the user last clicked (C10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week3")
$ws.Activate()

# DQ1 response 4 (row 9) took 30 minutes -> 0.020833333333333332 of a day
$ws.Range("C9").Value = [double](30) / 1440

# DQ2 response 2 (row 14) took 20 minutes -> 0.013888888888888888 of a day
$ws.Range("C14").Value = [double](20) / 1440

# Move the selection/active cell to where editing left off
$ws.Range("C10").Select()
